# Juno: check in to OLPRODLOC.
# Updates the French header labels on the "Contoso Chai Tea market trends 2023"
# worksheet and makes all header-row labels (besides "Date", already bold) bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# White color used by the header row text runs (RGB 255,255,255 -> 0xFFFFFF).
$headerWhite = 16777215

# --- B1: "Ventes totales de chaï (unités)" -> "Total des ventes de chaï (unités)" ---
$ws.Range("B1").Value = "Total des ventes de chaï (unités)"
$len = $ws.Range("B1").Text.Length
$ws.Range("B1").Characters(1, $len - 1).Font.Bold = $true
$ws.Range("B1").Characters(1, $len - 1).Font.Color = $headerWhite
$ws.Range("B1").Characters($len, 1).Font.Bold = $true
$ws.Range("B1").Characters($len, 1).Font.Color = $headerWhite

# --- C1: "Ventes de chaï artisanal (unités)" text unchanged, just becomes bold ---
$len = $ws.Range("C1").Text.Length
$ws.Range("C1").Characters(1, $len - 1).Font.Bold = $true
$ws.Range("C1").Characters($len, 1).Font.Bold = $true

# --- D1: "Ventes de chaï préparé (unités)" -> "Ventes de chaï préconfectionné (unités)" ---
$ws.Range("D1").Value = "Ventes de chaï préconfectionné (unités)"
$len = $ws.Range("D1").Text.Length
$ws.Range("D1").Characters(1, $len - 1).Font.Bold = $true
$ws.Range("D1").Characters(1, $len - 1).Font.Color = $headerWhite
$ws.Range("D1").Characters($len, 1).Font.Bold = $true
$ws.Range("D1").Characters($len, 1).Font.Color = $headerWhite

# --- E1: "Engagement sur les réseaux sociaux (vues)" text unchanged, just becomes bold ---
$len = $ws.Range("E1").Text.Length
$ws.Range("E1").Characters(1, $len - 1).Font.Bold = $true
$ws.Range("E1").Characters($len, 1).Font.Bold = $true

# --- F1: "Recherches en ligne du terme « chaï »" -> "Recherches en ligne de chaï" ---
$ws.Range("F1").Value = "Recherches en ligne de chaï"
$len = $ws.Range("F1").Text.Length
$ws.Range("F1").Characters(1, $len - 1).Font.Bold = $true
$ws.Range("F1").Characters(1, $len - 1).Font.Color = $headerWhite
$ws.Range("F1").Characters($len, 1).Font.Bold = $true
$ws.Range("F1").Characters($len, 1).Font.Color = $headerWhite
